$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the ambiguous column definition for column A ---
# (it previously spanned columns 1-2; re-asserting column B width narrows
# column As definition down to just column 1, as in the target file)
$ws.Columns("B:B").ColumnWidth = 60.7109375

# --- Insert a new row 13 to host the "Docentes responsaveis" value ---
# This shifts the previously mis-aligned content in rows 13-23 down to 14-24
$ws.Rows.Item(13).Insert()

# Clear the label cell that Insert() copied into A13 (new row 13 has no column-A label)
$ws.Range("A13").Clear()

# Give B13/C13 the same formatting as the other value cells before filling them in
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '2143261 - André Luis Ferraz'
$ws.Range("C13").Value = '2143261 - André Luis Ferraz'

# --- Row 10 (Objetivos): replace misplaced teacher name with the real objectives text ---
$ws.Range("B10").Value = 'Introdução à tecnologia de conversão de biomassa vegetal para estudantes de Engenharia Bioquímica, abordando os principais processos tecnológicos do setor e seus métodos de controle que incluem: celulose e papel; derivados de celulose; carvão vegetal e frações monoméricas por hidrólise.'
$ws.Range("C10").Value = 'Introdução à tecnologia de conversão de biomassa vegetal para estudantes de Engenharia Bioquímica, abordando os principais processos tecnológicos do setor e seus métodos de controle que incluem: celulose e papel; derivados de celulose; carvão vegetal e frações monoméricas por hidrólise.'

# --- Row 14 (Programa resumido): replace "Semestral" with the real Portuguese summary ---
$ws.Range("B14").Value = 'A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares.'
$ws.Range("C14").Value = 'A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares.'

# --- Row 15 (Short syllabus): unchanged text, left as-is ---

# --- Row 16 (Programa): replace stray date with the real Portuguese programme text ---
$ws.Range("B16").Value = '1. Breve introdução sobre a disponibilidade da biomassa: tipos de biomassa lignificada, produção por reflorestamento, resíduos agrícolas, características celulares, composição química.2. Produção de celuloses e papel: mercado mundial de celulose e papel, processos de polpação mecânica, kraft e sulfito; braqueamento de pastas celulósicas; recuperação de inorgânicos e geração de energia; métodos de controle de processo; características físico-químicas e métodos de produção de papel.3. Produção de derivados de celulose: formação do celulosato em meio alcalino, nitrato de celulose, xantato de celulose e a produção de fibras têxteis de "viscose", acetato de celulose, carboximetil celulose, etil e propilcelulose, alongamento da cadeia celulósica com epóxidos.4. Conversão térmica e produção de carvão vegetal: secagem da madeira e estabilização dimensional, processos termomecânicos e produção de aglomerados; energia de biomassa vegetal, queima para geração direta de energia; produção de carvão vegetal.5. Produção de açúcares e derivados por hidrólise: hidrólise ácida e processos de pré-tratamento para desestruturação da parede celular.6. Processos integrados para a conversão de biomassa: indústrias modernas que aplicam os conceitos de uso integrado da biomassa vegetal.'
$ws.Range("C16").Value = '1. Breve introdução sobre a disponibilidade da biomassa: tipos de biomassa lignificada, produção por reflorestamento, resíduos agrícolas, características celulares, composição química.2. Produção de celuloses e papel: mercado mundial de celulose e papel, processos de polpação mecânica, kraft e sulfito; braqueamento de pastas celulósicas; recuperação de inorgânicos e geração de energia; métodos de controle de processo; características físico-químicas e métodos de produção de papel.3. Produção de derivados de celulose: formação do celulosato em meio alcalino, nitrato de celulose, xantato de celulose e a produção de fibras têxteis de "viscose", acetato de celulose, carboximetil celulose, etil e propilcelulose, alongamento da cadeia celulósica com epóxidos.4. Conversão térmica e produção de carvão vegetal: secagem da madeira e estabilização dimensional, processos termomecânicos e produção de aglomerados; energia de biomassa vegetal, queima para geração direta de energia; produção de carvão vegetal.5. Produção de açúcares e derivados por hidrólise: hidrólise ácida e processos de pré-tratamento para desestruturação da parede celular.6. Processos integrados para a conversão de biomassa: indústrias modernas que aplicam os conceitos de uso integrado da biomassa vegetal.'

# --- Row 17 (Syllabus): unchanged text, left as-is ---

# --- Row 19 (Metodo): replace misplaced teacher name with the real evaluation method text ---
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# --- Row 20 (Criterio): replace with the final-grade formula text ---
$ws.Range("B20").Value = 'A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2'
$ws.Range("C20").Value = 'A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2'

# --- Row 21 (Norma de recuperacao): replace with the make-up exam text ---
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2'

# --- Row 22 (Bibliografia): add the reference list ---
$ws.Range("B22").Value = '1. EK, M., GELLERSTEDT, G., HENRIKSSON, G. Wood Chemistry and Wood Biotechnology (Vol 1) e Pulping Chemistry and Technology (Vol 2). Berlin: Walter de Gruyter, 2009.
2. KLEMM, D., PHILIPP, B., HEINZE, T., HEINZE, U., WAGENKNECHT, U. Comprehensive Cellulose Chemistry (Volume 2-Functionalization of Cellulose). Berlin: Wyley, 1998.
3. FENGEL, D., WEGENER, G. Wood Chemistry, Ultrastruture, Reactions. Berlin: Walter de Gruyter,1989.'
$ws.Range("C22").Value = '1. EK, M., GELLERSTEDT, G., HENRIKSSON, G. Wood Chemistry and Wood Biotechnology (Vol 1) e Pulping Chemistry and Technology (Vol 2). Berlin: Walter de Gruyter, 2009.
2. KLEMM, D., PHILIPP, B., HEINZE, T., HEINZE, U., WAGENKNECHT, U. Comprehensive Cellulose Chemistry (Volume 2-Functionalization of Cellulose). Berlin: Wyley, 1998.
3. FENGEL, D., WEGENER, G. Wood Chemistry, Ultrastruture, Reactions. Berlin: Walter de Gruyter,1989.'

# --- Row 24 (Requisitos value): unchanged text, left as-is ---

